$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 302
$ws.Range("I4").Value = 269.45456
$ws.Range("J4").Value = 391.5
$ws.Range("K4").Value = 269.45456
$ws.Range("L4").Value = 391.5
$ws.Range("M4").Value = -155.45456
$ws.Range("N4").Value = -619.5
$ws.Range("H6").Value = 368.8889
$ws.Range("I6").Value = 368.8889
$ws.Range("K6").Value = 1106.6667
$ws.Range("M6").Value = -994.6667
$ws.Range("H101").Value = 158672.2
$ws.Range("I101").Value = 1343.5
$ws.Range("K101").Value = 4030.5
$ws.Range("M101").Value = -2408.5
$ws.Range("H106").Value = 2050.4443
$ws.Range("I106").Value = 2422.1428
$ws.Range("J106").Value = 749.5
$ws.Range("K106").Value = 2422.1428
$ws.Range("L106").Value = 749.5
$ws.Range("M106").Value = -1791.1428
$ws.Range("N106").Value = -2011.5
$ws.Range("H125").Value = 7892
$ws.Range("I125").Value = 8004
$ws.Range("J125").Value = 7780
$ws.Range("K125").Value = 72036
$ws.Range("L125").Value = 70020
$ws.Range("M125").Value = -69576
$ws.Range("N125").Value = -74940
$ws.Range("H132").Value = 2050.5652
$ws.Range("I132").Value = 1689.3636
$ws.Range("K132").Value = 5068.0908
$ws.Range("M132").Value = -2538.0908
$ws.Range("H135").Value = 1525.619
$ws.Range("I135").Value = 1331.6111
$ws.Range("J135").Value = 2689.6667
$ws.Range("K135").Value = 11984.4999
$ws.Range("L135").Value = 24207.0003
$ws.Range("M135").Value = -9449.499900000001
$ws.Range("N135").Value = -29277.0003
$ws.Range("H137").Value = 416458.75
$ws.Range("I137").Value = 2334.6667
$ws.Range("K137").Value = 7004.000100000001
$ws.Range("M137").Value = -4454.000100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2070.2083
$ws.Range("I74").Value = 1504.6316
$ws.Range("K74").Value = 1504.6316
$ws.Range("M74").Value = -630.6315999999999
$ws.Range("H77").Value = 2070.2083
$ws.Range("I77").Value = 1504.6316
$ws.Range("K77").Value = 7523.157999999999
$ws.Range("M77").Value = -3155.157999999999
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H128").Value = 71100
$ws.Range("J128").Value = 71100
$ws.Range("L128").Value = 71100
$ws.Range("N128").Value = -81060
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H138").Value = 70595.75
$ws.Range("J138").Value = 69331
$ws.Range("L138").Value = 69331
$ws.Range("N138").Value = -79611
$ws.Range("N114").ClearContents()
$ws.Range("N133").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H94").Value = 3529.9524
$ws.Range("I94").Value = 3771.7058
$ws.Range("K94").Value = 3771.7058
$ws.Range("M94").Value = -3320.7058
$ws.Range("H105").Value = 103367.2
$ws.Range("I105").Value = 250899.75
$ws.Range("K105").Value = 250899.75
$ws.Range("M105").Value = -249152.75
$ws.Range("H134").Value = 3230.976
$ws.Range("J134").Value = 7700
$ws.Range("L134").Value = 23100
$ws.Range("N134").Value = -28170
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3428.4211
$ws.Range("I31").Value = 2108.1667
$ws.Range("J31").Value = 4037.7693
$ws.Range("K31").Value = 2108.1667
$ws.Range("L31").Value = 4037.7693
$ws.Range("M31").Value = -1813.1667
$ws.Range("N31").Value = -4627.7693
$ws.Range("H34").Value = 3428.4211
$ws.Range("I34").Value = 2108.1667
$ws.Range("J34").Value = 4037.7693
$ws.Range("K34").Value = 2108.1667
$ws.Range("L34").Value = 4037.7693
$ws.Range("M34").Value = -1906.1667
$ws.Range("N34").Value = -4441.7693
$ws.Range("H59").Value = 114998.5
$ws.Range("J59").Value = 114998.5
$ws.Range("L59").Value = 114998.5
$ws.Range("N59").Value = -117288.5
$ws.Range("H105").Value = 371836.66
$ws.Range("I105").Value = 371836.66
$ws.Range("K105").Value = 371836.66
$ws.Range("M105").Value = -370089.66
$ws.Range("H132").Value = 2656
$ws.Range("I132").Value = 2299.8
$ws.Range("K132").Value = 6899.400000000001
$ws.Range("M132").Value = -4369.400000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1435357
$ws.Range("I9").Value = 5014999.5
$ws.Range("J9").Value = 3500
$ws.Range("K9").Value = 15044998.5
$ws.Range("L9").Value = 10500
$ws.Range("M9").Value = -15044774.5
$ws.Range("N9").Value = -10948
$ws.Range("H17").Value = 550
$ws.Range("I17").Value = 550
$ws.Range("K17").Value = 1650
$ws.Range("M17").Value = -1481
$ws.Range("H23").Value = 77044.69500000001
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 83463.836
$ws.Range("K23").Value = 45
$ws.Range("L23").Value = 250391.508
$ws.Range("M23").Value = 190
$ws.Range("N23").Value = -250861.508
$ws.Range("H41").Value = 699
$ws.Range("I41").Value = 236.625
$ws.Range("K41").Value = 709.875
$ws.Range("M41").Value = -371.875
$ws.Range("H97").Value = 192
$ws.Range("J97").Value = 249
$ws.Range("L97").Value = 747
$ws.Range("N97").Value = -1739

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 1000000000
$ws.Range("J64").Value = 1000000000
$ws.Range("L64").Value = 1000000000
$ws.Range("N64").Value = -1000000496
$ws.Range("H67").Value = 1000000000
$ws.Range("J67").Value = 1000000000
$ws.Range("L67").Value = 1000000000
$ws.Range("N67").Value = -1000001716
$ws.Range("H102").Value = 2657.8572
$ws.Range("I102").Value = 2976.5
$ws.Range("J102").Value = 2233
$ws.Range("K102").Value = 2976.5
$ws.Range("L102").Value = 2233
$ws.Range("M102").Value = -1354.5
$ws.Range("N102").Value = -5477
$ws.Range("H103").Value = 48575.25
$ws.Range("J103").Value = 44767.332
$ws.Range("L103").Value = 44767.332
$ws.Range("N103").Value = -47111.332
$ws.Range("H126").Value = 4442.579
$ws.Range("J126").Value = 5338.231
$ws.Range("L126").Value = 16014.693
$ws.Range("N126").Value = -20954.693
$ws.Range("H132").Value = 4027.923
$ws.Range("I132").Value = 3152.3125
$ws.Range("K132").Value = 9456.9375
$ws.Range("M132").Value = -6926.9375
$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2570.5833
$ws.Range("I132").Value = 2156.75
$ws.Range("K132").Value = 6470.25
$ws.Range("M132").Value = -3940.25
$ws.Range("H136").Value = 8714.4
$ws.Range("I136").Value = 18787
$ws.Range("J136").Value = 1999.3334
$ws.Range("K136").Value = 56361
$ws.Range("L136").Value = 5998.0002
$ws.Range("M136").Value = -53811
$ws.Range("N136").Value = -11098.0002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H123").Value = 49999.5
$ws.Range("J123").Value = 49999
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -59799
$ws.Range("H136").Value = 751.4
$ws.Range("I136").Value = 473.77777
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 1421.33331
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = 1128.66669
$ws.Range("N136").Value = -14850
$ws.Range("N109").ClearContents()
